$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row labels: "_old" suffix -> "_FV2310", "_new" suffix -> "_FV2404"
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $baseNames[$i] + "_FV2310"
    $ws.Range($newCols[$i] + "1").Value = $baseNames[$i] + "_FV2404"
}

# Turn the data range into an Excel Table ("Table1") with an AutoFilter
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U92"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium2"
$tbl.ShowTableStyleRowStripes = $true
$tbl.ShowTableStyleFirstColumn = $false
$tbl.ShowTableStyleLastColumn = $false
$tbl.ShowTableStyleColumnStripes = $false

# Freeze the header row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
